# Actualizar funcionalidades de salud financiera
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: CASA LEO FERRETERIA's rol changes from "cliente_premium" to "admin"
$ws.Range("D4").Value = "admin"

# New row 7: RAFAEL FERRETERIA
$ws.Range("A7").Value = 20246
$ws.Range("B7").Value = "RAFAEL FERRETERIA"
$ws.Range("C7").Value = 20246
$ws.Range("D7").Value = "vendedor_estandar "
